# [master] added last lesson
#
# Slide 4, "Content Placeholder 7" shape: the paragraph that reads
# "consume per min (l/min). " gets the closing paren pulled out of the
# first run into its own run, and the paragraph is split in two right
# after the ". " run (leaving a new, otherwise-empty paragraph that
# reuses the original endParaRPr).

$p = $ppt.ActivePresentation

# Find the slide + shape that holds the "consume per min" text instead of
# hard-coding indices, so the script is resilient to minor reshuffles.
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            if ($shape.TextFrame.TextRange.Text -like "*consume per min*") {
                $targetSlide = $slide
                $targetShape = $shape
                break
            }
        }
    }
    if ($targetShape -ne $null) { break }
}

$tr = $targetShape.TextFrame.TextRange

# Locate the paragraph containing "consume per min" dynamically.
$paraIndex = -1
for ($i = 1; $i -le $tr.Paragraphs(0, -1).Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text -like "*consume per min*") {
        $paraIndex = $i
        break
    }
}

$para = $tr.Paragraphs($paraIndex, 1)

# Run 1 holds "consume per min (l/min)"; split the trailing ")" into its
# own run by nudging its (already-correct) bold state, which forces the
# engine to break the run at that boundary.
$run1 = $para.Runs(1, 1)
$run1Text = $run1.Text
$closeParenOffset = $run1Text.Length - 1  # 0-based offset of ")" within run1
$closeParenPos = $run1.Start + $closeParenOffset  # 1-based absolute position

$tr2 = $targetShape.TextFrame.TextRange
$closeParenRange = $tr2.Characters($closeParenPos, 1)
$closeParenRange.Font.Bold = $false

# Re-fetch the paragraph/runs (now 3 runs: "...l/min", ")", ". ").
$tr3 = $targetShape.TextFrame.TextRange
$para2 = $tr3.Paragraphs($paraIndex, 1)
$periodRun = $para2.Runs(3, 1)

# Split the paragraph in two right after the ". " run.
$periodRun.InsertAfter([char]13)

# The split leaves a stray empty run at the start of the freshly created
# paragraph (before its endParaRPr); clearing its text removes the run
# entirely so the new paragraph is just pPr + endParaRPr, matching a
# plain Enter keypress at the end of the line.
$tr4 = $targetShape.TextFrame.TextRange
$newPara = $tr4.Paragraphs($paraIndex + 1, 1)
$strayRun = $newPara.Runs(1, 1)
$strayRun.Text = ""
